# Update the dSF (column F) values for specific rows as per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 5
    10 = 0
    15 = -8
    16 = -2
    17 = -12
    18 = -6
    19 = -7
    20 = 0
    23 = -3
    26 = -8
    27 = 0
    29 = 9
    35 = -6
    38 = 5
    40 = -1
    43 = -4
    46 = -5
    47 = -6
    53 = -4
    57 = -4
    61 = 0
    63 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
